$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 5 new rows starting at row 24 (shifts old rows 25-29 down to 30-34) ---
$ws.Rows("24:28").Insert()

# --- Row 23: fill in newly-added cells (existing C23/D23/E23/I23 stay as-is) ---
$dateFmt = 'm/d/yy\ h:mm;@'
$ws.Cells.Item(23, 1).Value = 5
$ws.Cells.Item(23, 2).Value = "מעשי"
$ws.Cells.Item(23, 6).Value = 43999.537499999999
$ws.Cells.Item(23, 6).NumberFormat = $dateFmt
$ws.Cells.Item(23, 7).Value = 0.5
$ws.Cells.Item(23, 8).Value = 0.5

# --- Row 24: new row (chapter 6, SVM) ---
$ws.Cells.Item(24, 1).Value = 6
$ws.Cells.Item(24, 2).Value = "תאורטי"
$ws.Cells.Item(24, 3).Value = "קריאה+שאלות"
$ws.Cells.Item(24, 4).Value = "ISLR 337-359`nשאלות 3,6 ועוד כמה מההוראות בדרייב"
$ws.Cells.Item(24, 4).WrapText = $true
$ws.Cells.Item(24, 5).Value = 43999.59097222222
$ws.Cells.Item(24, 5).NumberFormat = $dateFmt
$ws.Cells.Item(24, 6).Value = 44000.605555555558
$ws.Cells.Item(24, 6).NumberFormat = $dateFmt
$ws.Cells.Item(24, 7).Value = 1
$ws.Cells.Item(24, 8).Value = 0.75
$ws.Cells.Item(24, 9).Value = "תרגיל 6 ב ISLR מרגיש מיותר ומנסה להעביר נקודה שדי פשוט להבין מהקריאה עצמה`nב Part 4 ב drive יוצא (לי, אולי אני בעייתי) ששני המודלים מסווגים הכל כ-0 בגלל שהדאטה לא מאוזן."
$ws.Cells.Item(24, 9).WrapText = $true
$ws.Rows(24).RowHeight = 45

# --- Row 25: new row (chapter 6, Bayesian learning) ---
$ws.Cells.Item(25, 1).Value = 6
$ws.Cells.Item(25, 2).Value = "תאורטי"
$ws.Cells.Item(25, 3).Value = "קריאה"
$ws.Cells.Item(25, 4).Value = "deep learning 5.6`nBishop 3.3-3.3.1`nחלק כלשהו מהקורס ב MIT`nIntoduction to Probability 8"
$ws.Cells.Item(25, 4).WrapText = $true
$ws.Cells.Item(25, 5).Value = 44000.615972222222
$ws.Cells.Item(25, 5).NumberFormat = $dateFmt
$ws.Cells.Item(25, 6).Value = 44004.611111111109
$ws.Cells.Item(25, 6).NumberFormat = $dateFmt
$ws.Cells.Item(25, 7).Value = 2.25
$ws.Cells.Item(25, 8).Value = 1.75
$ws.Cells.Item(25, 9).Value = "עדיף להתחיל מקריאה של הבלוג פוסטים לפני שנכנסים לקריאה בספרים (שרובה הרגישה לא הכי חיונית). `nאם אפשר כדאי (לדעתי) להוסיף תרגיל לפרק (גם אם משהו קטן של שימוש ב PYMC3)"
$ws.Cells.Item(25, 9).WrapText = $true
$ws.Rows(25).RowHeight = 60

# --- Remove the leftover blank rows (26-29) created by the Insert above ---
$ws.Range("E26:E29").Clear()

# --- Update view: scroll position + active selection ---
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I26").Select()
